$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 6.63718292801021
$ws.Range("D2").Value = 8.93459740334325
$ws.Range("E2").Value = 10.96374181033933
$ws.Range("F2").Value = 44.44119803542744
$ws.Range("G2").Value = 3.743964081058378
$ws.Range("I2").Value = 37.95340611643464
$ws.Range("K2").Value = 19.57135239643369
$ws.Range("L2").Value = 9.424842787455148
$ws.Range("M2").Value = 19.86435596176544

$ws.Range("C3").Value = 6.63946683161606
$ws.Range("D3").Value = 8.962783197479377
$ws.Range("E3").Value = 10.98157163812049
$ws.Range("F3").Value = 44.27765426140905
$ws.Range("G3").Value = 3.748030053671211
$ws.Range("I3").Value = 37.89105351662516
$ws.Range("K3").Value = 19.30215532147275
$ws.Range("L3").Value = 9.445332836830747
$ws.Range("M3").Value = 19.77687252305352

$ws.Range("C4").Value = 6.641812053616699
$ws.Range("D4").Value = 8.980963041199407
$ws.Range("E4").Value = 10.99402296886452
$ws.Range("F4").Value = 44.18877103296594
$ws.Range("G4").Value = 3.750654125121094
$ws.Range("I4").Value = 37.86180507877312
$ws.Range("K4").Value = 19.14067235887818
$ws.Range("L4").Value = 9.458892998494637
$ws.Range("M4").Value = 19.72769604185668

$ws.Range("C5").Value = 6.64300574626817
$ws.Range("D5").Value = 8.988591245522064
$ws.Range("E5").Value = 10.99947504515762
$ws.Range("F5").Value = 44.15546271921467
$ws.Range("G5").Value = 3.751755656151877
$ws.Range("I5").Value = 37.85215685432532
$ws.Range("K5").Value = 19.07590767195866
$ws.Range("L5").Value = 9.4646652782489
$ws.Range("M5").Value = 19.70881269765286

$ws.Range("C6").Value = 6.643218357609397
$ws.Range("D6").Value = 8.989871179255703
$ws.Range("E6").Value = 11.00040318942128
$ws.Range("F6").Value = 44.15010807430304
$ws.Range("G6").Value = 3.751940513124494
$ws.Range("I6").Value = 37.85069183815769
$ws.Range("K6").Value = 19.06521883349383
$ws.Range("L6").Value = 9.465638648393096
$ws.Range("M6").Value = 19.70574739451326

$ws.Range("C7").Value = 6.641827187380327
$ws.Range("D7").Value = 8.981065027907247
$ws.Range("E7").Value = 10.99409496684529
$ws.Range("F7").Value = 44.18831002064651
$ws.Range("G7").Value = 3.750668850207452
$ws.Range("I7").Value = 37.86166576896202
$ws.Range("K7").Value = 19.13979459538961
$ws.Range("L7").Value = 9.458969847566706
$ws.Range("M7").Value = 19.72743667310725

$ws.Range("C8").Value = 6.637775271492904
$ws.Range("D8").Value = 8.944134589024081
$ws.Range("E8").Value = 10.96957739729299
$ws.Range("F8").Value = 44.38242288410603
$ws.Range("G8").Value = 3.74533963290556
$ws.Range("I8").Value = 37.93003113780981
$ws.Range("K8").Value = 19.47779449731449
$ws.Range("L8").Value = 9.431704633236068
$ws.Range("M8").Value = 19.83325921943751

$ws.Range("C9").Value = 6.637270174457862
$ws.Range("D9").Value = 8.878641742224556
$ws.Range("E9").Value = 10.93343301788585
$ws.Range("F9").Value = 44.8538789208999
$ws.Range("G9").Value = 3.735895259676041
$ws.Range("I9").Value = 38.13577796104554
$ws.Range("K9").Value = 20.16703285052232
$ws.Range("L9").Value = 9.385999174587853
$ws.Range("M9").Value = 20.07605527777284

$ws.Range("C10").Value = 6.641378392020973
$ws.Range("D10").Value = 8.834740471373054
$ws.Range("E10").Value = 10.91415704913797
$ws.Range("F10").Value = 45.25441397062373
$ws.Range("G10").Value = 3.729561722364629
$ws.Range("I10").Value = 38.33049761128522
$ws.Range("K10").Value = 20.68442605844898
$ws.Range("L10").Value = 9.357140411152395
$ws.Range("M10").Value = 20.27486856038291

$ws.Range("C11").Value = 6.644206839723002
$ws.Range("D11").Value = 8.815682921774979
$ws.Range("E11").Value = 10.90696879469558
$ws.Range("F11").Value = 45.44805418081517
$ws.Range("G11").Value = 3.726810107180716
$ws.Range("I11").Value = 38.42846568634775
$ws.Range("K11").Value = 20.92111635410299
$ws.Range("L11").Value = 9.345034781378235
$ws.Range("M11").Value = 20.36949439913785

$ws.Range("C12").Value = 6.645414674964299
$ws.Range("D12").Value = 8.808597643772686
$ws.Range("E12").Value = 10.90447400796042
$ws.Range("F12").Value = 45.52299118333248
$ws.Range("G12").Value = 3.725786635440375
$ws.Range("I12").Value = 38.46690401712763
$ws.Range("K12").Value = 21.01084735209021
$ws.Range("L12").Value = 9.340597558882168
$ws.Range("M12").Value = 20.40590515222668

$ws.Range("C13").Value = 6.645148483550223
$ws.Range("D13").Value = 8.810117739819777
$ws.Range("E13").Value = 10.90500119869376
$ws.Range("F13").Value = 45.50678115572296
$ws.Range("G13").Value = 3.72600623732017
$ws.Range("I13").Value = 38.45856622435023
$ws.Range("K13").Value = 20.99151905594603
$ws.Range("L13").Value = 9.341546661498521
$ws.Range("M13").Value = 20.39803812689128

$ws.Range("C14").Value = 6.644303475327554
$ws.Range("D14").Value = 8.815097377773657
$ws.Range("E14").Value = 10.9067589929327
$ws.Range("F14").Value = 45.45418724337923
$ws.Range("G14").Value = 3.726725535375121
$ws.Range("I14").Value = 38.43160124209545
$ws.Range("K14").Value = 20.92849704516668
$ws.Range("L14").Value = 9.344666784531636
$ws.Range("M14").Value = 20.37247848923926

$ws.Range("C15").Value = 6.643803656715039
$ws.Range("D15").Value = 8.818164664415564
$ws.Range("E15").Value = 10.90786528593751
$ws.Range("F15").Value = 45.42218049352984
$ws.Range("G15").Value = 3.727168532486491
$ws.Range("I15").Value = 38.41525857496669
$ws.Range("K15").Value = 20.88990485788295
$ws.Range("L15").Value = 9.346597080029619
$ws.Range("M15").Value = 20.35689702422881

$ws.Range("C16").Value = 6.641212752390828
$ws.Range("D16").Value = 8.836004311137216
$ws.Range("E16").Value = 10.91465862073965
$ws.Range("F16").Value = 45.24198633881298
$ws.Range("G16").Value = 3.729744143108158
$ws.Range("I16").Value = 38.32428327545903
$ws.Range("K16").Value = 20.66897680030866
$ws.Range("L16").Value = 9.357952106665818
$ws.Range("M16").Value = 20.26876673854555

$ws.Range("C17").Value = 6.639868238393952
$ws.Range("D17").Value = 8.847182328882926
$ws.Range("E17").Value = 10.91923090653081
$ws.Range("F17").Value = 45.13434819958271
$ws.Range("G17").Value = 3.731357288262729
$ws.Range("I17").Value = 38.27087105325328
$ws.Range("K17").Value = 20.53372005697921
$ws.Range("L17").Value = 9.365179829128884
$ws.Range("M17").Value = 20.21575688643551

$ws.Range("C18").Value = 6.639185308374033
$ws.Range("D18").Value = 8.853697593961108
$ws.Range("E18").Value = 10.92200953033939
$ws.Range("F18").Value = 45.07351618064823
$ws.Range("G18").Value = 3.732297327189445
$ws.Range("I18").Value = 38.24103424760985
$ws.Range("K18").Value = 20.45605393679443
$ws.Range("L18").Value = 9.369433259557413
$ws.Range("M18").Value = 20.18566239310076

$ws.Range("C19").Value = 6.638969642528806
$ws.Range("D19").Value = 8.855918313668584
$ws.Range("E19").Value = 10.92297587536993
$ws.Range("F19").Value = 45.05310575434489
$ws.Range("G19").Value = 3.732617707413121
$ws.Range("I19").Value = 38.23108420719533
$ws.Range("K19").Value = 20.42978259584944
$ws.Range("L19").Value = 9.370889929624688
$ws.Range("M19").Value = 20.17554152713007

$ws.Range("C20").Value = 6.640002016967958
$ws.Range("D20").Value = 8.845983512372895
$ws.Range("E20").Value = 10.91872878223752
$ws.Range("F20").Value = 45.14569507359678
$ws.Range("G20").Value = 3.731184304391835
$ws.Range("I20").Value = 38.27646540133798
$ws.Range("K20").Value = 20.54810556777929
$ws.Range("L20").Value = 9.36440046643842
$ws.Range("M20").Value = 20.22135912509512

$ws.Range("C21").Value = 6.644547973047684
$ws.Range("D21").Value = 8.8136311705903
$ws.Range("E21").Value = 10.90623651903436
$ws.Range("F21").Value = 45.4695919571587
$ws.Range("G21").Value = 3.72651375885199
$ws.Range("I21").Value = 38.43948523662365
$ws.Range("K21").Value = 20.94700608126089
$ws.Range("L21").Value = 9.343746342248922
$ws.Range("M21").Value = 20.37997048175599

$ws.Range("C22").Value = 6.648315572852133
$ws.Range("D22").Value = 8.79325299992292
$ws.Range("E22").Value = 10.89939664892947
$ws.Range("F22").Value = 45.69064126264001
$ws.Range("G22").Value = 3.723569093081502
$ws.Range("I22").Value = 38.55383220108657
$ws.Range("K22").Value = 21.2082626993789
$ws.Range("L22").Value = 9.331103978619987
$ws.Range("M22").Value = 20.48698952232531

$ws.Range("C23").Value = 6.646232268310549
$ws.Range("D23").Value = 8.804059103132206
$ws.Range("E23").Value = 10.90292603704196
$ws.Range("F23").Value = 45.57181857154791
$ws.Range("G23").Value = 3.725130892903272
$ws.Range("I23").Value = 38.49209289356995
$ws.Range("K23").Value = 21.06880344997725
$ws.Range("L23").Value = 9.337773125919377
$ws.Range("M23").Value = 20.42957234335784

$ws.Range("C24").Value = 6.639941255100047
$ws.Range("D24").Value = 8.846525220499316
$ws.Range("E24").Value = 10.91895532555046
$ws.Range("F24").Value = 45.14056187694067
$ws.Range("G24").Value = 3.731262471081279
$ws.Range("I24").Value = 38.27393348411848
$ws.Range("K24").Value = 20.5416015763359
$ws.Range("L24").Value = 9.364752510692174
$ws.Range("M24").Value = 20.21882516268166

$ws.Range("C25").Value = 6.636615878731263
$ws.Range("D25").Value = 8.895618414182719
$ws.Range("E25").Value = 10.94193284621294
$ws.Range("F25").Value = 44.71672529730063
$ws.Range("G25").Value = 3.738343337131752
$ws.Range("I25").Value = 38.07245026655443
$ws.Range("K25").Value = 19.97827338936807
$ws.Range("L25").Value = 9.397533896434144
$ws.Range("M25").Value = 20.00670145524141
